$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.063.70'
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").Value = '1.834.59'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6335'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07553'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.93'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07744'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.69%  '

$ws.Range("D12").Value = '1.830.40'
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.009'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6720'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.36'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009640'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.090'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.21%  '

$ws.Range("D18").Value = '29.099.75'
$ws.Range("E18").Value = '  +0.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.60'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.87'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.203'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.65'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("E25").Value = '  +3.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.551'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.90%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.95'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.96%  '

$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.128'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.080'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.36%  '

$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05383'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.43%  '

$ws.Range("E33").Value = '  +2.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7469'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.142'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.48%  '

$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").Value = '1.245.22'
$ws.Range("E37").Value = '  -2.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.761'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.45%  '

$ws.Range("E39").Value = '  +0.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.638'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9066'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.01'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.54%  '

$ws.Range("D44").Value = '1.981.84'
$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.99'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.75%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5117'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4095'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.122'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.655'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.782'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.21%  '
